$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maçlar")

# Add missing scores for row 10 (Ahmet Minguzzi Grubu: Ajans Of - Ravager)
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 24

# Add missing scores for row 11 (Eren Bülbül Grubu: Araklı 1961 Spor - Hubuş FK)
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3

# Swap the match data between row 12 and row 13 (date/time stayed per row, but
# the actual match details moved from one slot to the other)
$c12 = $ws.Range("C12").Value()
$d12 = $ws.Range("D12").Value()
$e12 = $ws.Range("E12").Value()
$b12 = $ws.Range("B12").Value()

$c13 = $ws.Range("C13").Value()
$d13 = $ws.Range("D13").Value()
$e13 = $ws.Range("E13").Value()
$b13 = $ws.Range("B13").Value()

$ws.Range("B12").Value = $b13
$ws.Range("C12").Value = $c13
$ws.Range("D12").Value = $d13
$ws.Range("E12").Value = $e13

$ws.Range("B13").Value = $b12
$ws.Range("C13").Value = $c12
$ws.Range("D13").Value = $d12
$ws.Range("E13").Value = $e12

# Update view state to match the saved selection (scroll position follows in Excel's UI)
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I16").Select()
